$d = $word.ActiveDocument

# 1) Merge the ": " run with the following "Proyek yang akan dikerjakan..." run
#    into a single run (same text, same formatting) under "Deskripsi Proyek".
$d.Content.Find.Execute(
    ": Proyek yang akan dikerjakan adalah membuat aplikasi untuk menghubungkan Game Designer dengan Game Publisher untuk membuat jadwal pertemuan yang akan dilakukan dan akan memberi notifikasi 1 hari sebelum hari pertemuan. Nantinya didalam aplikasi ini masing-masing user/pengguna dapat me-review hasil pertemuan yang telah dilakukan",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Proyek yang akan dikerjakan adalah membuat aplikasi untuk menghubungkan Game Designer dengan Game Publisher untuk membuat jadwal pertemuan yang akan dilakukan dan akan memberi notifikasi 1 hari sebelum hari pertemuan. Nantinya didalam aplikasi ini masing-masing user/pengguna dapat me-review hasil pertemuan yang telah dilakukan",
    2)

# 2) Merge the ": " run with the following Stakeholder list run into a single run.
$d.Content.Find.Execute(
    ": Ketua Asosiasi Game Developer Semarang(Pemilik proyek), Wakil Ketua Asosiasi Game Developer Semarang (Penanggung jawab pengoprasian), Akbar Karunia Octaviantono(Ketua tim), Dani Rahman Hakim(Anggota), Ludi Agustia Aryanto(Anggota), Reganda Dhynar An Nura(Anggota), Andika Wahyu Rama Ardiansyah(Anggota)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Ketua Asosiasi Game Developer Semarang(Pemilik proyek), Wakil Ketua Asosiasi Game Developer Semarang (Penanggung jawab pengoprasian), Akbar Karunia Octaviantono(Ketua tim), Dani Rahman Hakim(Anggota), Ludi Agustia Aryanto(Anggota), Reganda Dhynar An Nura(Anggota), Andika Wahyu Rama Ardiansyah(Anggota)",
    2)

# 3) Extend the "Reminder." bullet under Deskripsi Fungsionalitas with more detail.
$d.Content.Find.Execute(
    "Reminder.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Reminder appoinment jika jadwal sudah dekat.",
    2)
